$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.22896466666667
$ws.Range("N2").Value = 93.686894
$ws.Range("O2").Value = 0.2877106972998646
$ws.Range("P2").Value = 0.2877106972998646
$ws.Range("Q2").Value = 1272.686217995506
$ws.Range("R2").Value = 11454.17596195956
$ws.Range("S2").Value = 0.006118273369924755
$ws.Range("T2").Value = 0.006118273369924755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3726247238124506
$ws.Range("P3").Value = 0.3726247238124505
$ws.Range("Q3").Value = 1648.302808797617
$ws.Range("R3").Value = 14834.72527917855
$ws.Range("S3").Value = 0.007924001248730616
$ws.Range("T3").Value = 0.007924001248730614

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.36964133333333
$ws.Range("N4").Value = 76.108924
$ws.Range("O4").Value = 0.2337290805561598
$ws.Range("P4").Value = 0.2337290805561598
$ws.Range("Q4").Value = 1033.898921243642
$ws.Range("R4").Value = 9305.090291192777
$ws.Range("S4").Value = 0.004970334515762974
$ws.Range("T4").Value = 0.004970334515762974

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.49855033333333
$ws.Range("N5").Value = 34.495651
$ws.Range("O5").Value = 0.1059354983315251
$ws.Range("P5").Value = 0.1059354983315251
$ws.Range("Q5").Value = 468.6049215003638
$ws.Range("R5").Value = 4217.444293503275
$ws.Range("S5").Value = 0.002252757177450224
$ws.Range("T5").Value = 0.002252757177450224

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.22896466666667
$ws.Range("N6").Value = 93.686894
$ws.Range("O6").Value = 0.2877106972998646
$ws.Range("P6").Value = 0.2877106972998646
$ws.Range("Q6").Value = 52754.75606967116
$ws.Range("R6").Value = 474792.8046270404
$ws.Range("S6").Value = 0.2536116244790553
$ws.Range("T6").Value = 0.2536116244790553

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3726247238124506
$ws.Range("P7").Value = 0.3726247238124505
$ws.Range("Q7").Value = 68324.62815856401
$ws.Range("R7").Value = 614921.6534270761
$ws.Range("S7").Value = 0.3284617583358079
$ws.Range("T7").Value = 0.3284617583358079

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 25.36964133333333
$ws.Range("N8").Value = 76.108924
$ws.Range("O8").Value = 0.2337290805561598
$ws.Range("P8").Value = 0.2337290805561598
$ws.Range("Q8").Value = 42856.66381836867
$ws.Range("R8").Value = 385709.9743653181
$ws.Range("S8").Value = 0.2060278340852346
$ws.Range("T8").Value = 0.2060278340852346

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.49855033333333
$ws.Range("N9").Value = 34.495651
$ws.Range("O9").Value = 0.1059354983315251
$ws.Range("P9").Value = 0.1059354983315251
$ws.Range("Q9").Value = 19424.37812026844
$ws.Range("R9").Value = 174819.4030824159
$ws.Range("S9").Value = 0.09338016999018613
$ws.Range("T9").Value = 0.09338016999018613

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.22896466666667
$ws.Range("N10").Value = 93.686894
$ws.Range("O10").Value = 0.2877106972998646
$ws.Range("P10").Value = 0.2877106972998646
$ws.Range("Q10").Value = 3153.044982023451
$ws.Range("R10").Value = 28377.40483821106
$ws.Range("S10").Value = 0.01515785342444643
$ws.Range("T10").Value = 0.01515785342444643

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3726247238124506
$ws.Range("P11").Value = 0.3726247238124505
$ws.Range("Q11").Value = 4083.624719626559
$ws.Range("R11").Value = 36752.62247663904
$ws.Range("S11").Value = 0.01963149441046736
$ws.Range("T11").Value = 0.01963149441046736

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 25.36964133333333
$ws.Range("N12").Value = 76.108924
$ws.Range("O12").Value = 0.2337290805561598
$ws.Range("P12").Value = 0.2337290805561598
$ws.Range("Q12").Value = 2561.456044272363
$ws.Range("R12").Value = 23053.10439845127
$ws.Range("S12").Value = 0.01231386659359561
$ws.Range("T12").Value = 0.01231386659359561

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.49855033333333
$ws.Range("N13").Value = 34.495651
$ws.Range("O13").Value = 0.1059354983315251
$ws.Range("P13").Value = 0.1059354983315251
$ws.Range("Q13").Value = 1160.955760655084
$ws.Range("R13").Value = 10448.60184589576
$ws.Range("S13").Value = 0.005581143736485264
$ws.Range("T13").Value = 0.005581143736485264

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.22896466666667
$ws.Range("N14").Value = 93.686894
$ws.Range("O14").Value = 0.2877106972998646
$ws.Range("P14").Value = 0.2877106972998646
$ws.Range("Q14").Value = 2667.351668555595
$ws.Range("R14").Value = 24006.16501700035
$ws.Range("S14").Value = 0.01282294602643813
$ws.Range("T14").Value = 0.01282294602643813

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3726247238124506
$ws.Range("P15").Value = 0.3726247238124505
$ws.Range("Q15").Value = 3454.585415606913
$ws.Range("R15").Value = 31091.26874046222
$ws.Range("S15").Value = 0.01660746981744471
$ws.Range("T15").Value = 0.0166074698174447

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 25.36964133333333
$ws.Range("N16").Value = 76.108924
$ws.Range("O16").Value = 0.2337290805561598
$ws.Range("P16").Value = 0.2337290805561598
$ws.Range("Q16").Value = 2166.89076514129
$ws.Range("R16").Value = 19502.01688627161
$ws.Range("S16").Value = 0.0104170453615666
$ws.Range("T16").Value = 0.0104170453615666

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.49855033333333
$ws.Range("N17").Value = 34.495651
$ws.Range("O17").Value = 0.1059354983315251
$ws.Range("P17").Value = 0.1059354983315251
$ws.Range("Q17").Value = 982.1227743206159
$ws.Range("R17").Value = 8839.104968885544
$ws.Range("S17").Value = 0.004721427427403525
$ws.Range("T17").Value = 0.004721427427403525

